$d = $word.ActiveDocument

# Locate the paragraph holding the whole "Convention n ... {% endif %}" Jinja
# block (the opening title line of the fiche) without assuming its absolute
# paragraph index.
$locate = $d.Content
$locate.Find.ClearFormatting()
$null = $locate.Find.Execute("Convention n", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$p1 = $d.Range($locate.Start, $locate.Start).Paragraphs(1)
$r = $p1.Range

# Match the full existing paragraph text (without its trailing paragraph mark)
$oldText = "Convention n" + [char]176 + " {% if convention.numero %}{{convention.numero}}{% else %}Le num" + [char]233 + "ro de la convention sera d" + [char]233 + "fini et ajout" + [char]233 + " ici une fois la convention valid" + [char]233 + "e{% endif %}"

# New content: 6 paragraphs (5 with text + a trailing blank one), using ^p to split
$newText = "{% if convention.is_avenant() -%}^p" + `
           "Avenant n" + [char]176 + " {{ convention.numero }} " + [char]224 + " la convention n" + [char]176 + " {{ convention.parent.numero }}^p" + `
           "{%- else -%}^p" + `
           "Convention n" + [char]176 + " {{ convention.numero }}^p" + `
           "{%- endif %}^p"

$startPos = $r.Start
$ok = $r.Find.Execute($oldText, $false, $false, $false, $false, $false, $true, 1, $false, $newText, 2)
Write-Output "replace-ok:$ok"

Write-Output $d.Paragraphs.Count

# Make sure the formatting across the newly created paragraphs is uniform bold,
# matching the surrounding run formatting (the removed run was italic/non-bold).
$pNew = $d.Range($startPos, $startPos).Paragraphs(1)
for ($i = 0; $i -lt 5; $i++) {
    $pNew = $pNew.Next()
}
$endRange = $d.Range($startPos, $pNew.Range.End)
$endRange.Font.Bold = 1

# Further down, the empty paragraph right before "Nombre de logements :" gets
# an explicit en-US paragraph-mark language tag.
$rng = $d.Content
$rng.Find.ClearFormatting()
$null = $rng.Find.Execute("Nombre de logements", $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
$nbLogementsPara = $d.Range($rng.Start, $rng.Start).Paragraphs(1)
$blankPara = $nbLogementsPara.Previous()
$blankPara.Range.LanguageID = "en-US"

Write-Output "done"
